# The document is edited so that every paragraph's <w:contextualSpacing .../>
# element (always w:val="0") is removed from its <w:pPr>. There is no
# ContextualSpacing property exposed on Word's ParagraphFormat object model
# here, so we drop down to the underlying OOXML package via
# Document.WordOpenXML, strip the element, and write the package back.

$d = $word.ActiveDocument

$xml = $d.WordOpenXML

$before = ([regex]::Matches($xml, '<w:contextualSpacing\b[^/]*/>')).Count

$xml = [regex]::Replace($xml, '<w:contextualSpacing\b[^/]*/>', '')

$after = ([regex]::Matches($xml, '<w:contextualSpacing\b[^/]*/>')).Count

$d.WordOpenXML = $xml

Write-Host "contextualSpacing occurrences before:" $before "after:" $after
